$wb = $excel.ActiveWorkbook

# Add the new worksheet after the existing one (will become "Sheet1", second tab)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws2.Name = "Sheet1"

# Header row (order chosen to match shared-string insertion order)
$ws2.Range("A1").Value = "Target"
$ws2.Range("C1").Value = "Momentum"
$ws2.Range("B1").Value = "l1"
$ws2.Range("D1").Value = "Attempt"
$ws2.Range("E1").Value = "Epoc"

# Data rows
$ws2.Range("A2").Value = 20
$ws2.Range("B2").Value = 10
$ws2.Range("C2").Value = 0.1
$ws2.Range("D2").Value = 1
$ws2.Range("E2").Value = 2964

$ws2.Range("A3").Value = 20
$ws2.Range("B3").Value = 10
$ws2.Range("C3").Value = 0.1
$ws2.Range("D3").Value = 1
$ws2.Range("E3").Value = 1248

$ws2.Range("A4").Value = 20
$ws2.Range("B4").Value = 10
$ws2.Range("C4").Value = 0.1
$ws2.Range("D4").Value = 1
$ws2.Range("E4").Value = 299

$ws2.Range("A5").Value = 20
$ws2.Range("B5").Value = 10
$ws2.Range("C5").Value = 0
$ws2.Range("D5").Value = 1
$ws2.Range("E5").Value = 3055

$ws2.Columns("C").AutoFit() | Out-Null

$ws2.Range("C6").Select()

$ws2.Activate()
